$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Cells.Item(2, 4) "25.720.22"
Set-TextValue $ws.Cells.Item(2, 5) "  -4.07%  "

# Row 3
Set-TextValue $ws.Cells.Item(3, 4) "1.817.34"
Set-TextValue $ws.Cells.Item(3, 5) "  -2.87%  "

# Row 4
Set-TextValue $ws.Cells.Item(4, 5) "  +0.07%  "

# Row 5
Set-TextValue $ws.Cells.Item(5, 4) "278.99"
Set-TextValue $ws.Cells.Item(5, 5) "  -7.40%  "

# Row 6
Set-TextValue $ws.Cells.Item(6, 5) "  +0.06%  "

# Row 7
Set-TextValue $ws.Cells.Item(7, 4) "0.5087"
Set-TextValue $ws.Cells.Item(7, 5) "  -4.74%  "

# Row 8
Set-TextValue $ws.Cells.Item(8, 4) "0.3537"
Set-TextValue $ws.Cells.Item(8, 5) "  -5.53%  "

# Row 9
Set-TextValue $ws.Cells.Item(9, 4) "44.35"
Set-TextValue $ws.Cells.Item(9, 5) "  -2.46%  "

# Row 10
Set-TextValue $ws.Cells.Item(10, 4) "0.06686"
Set-TextValue $ws.Cells.Item(10, 5) "  -7.03%  "

# Row 11
Set-TextValue $ws.Cells.Item(11, 4) "19.91"
Set-TextValue $ws.Cells.Item(11, 5) "  -8.05%  "

# Row 12
Set-TextValue $ws.Cells.Item(12, 4) "0.8242"
Set-TextValue $ws.Cells.Item(12, 5) "  -7.32%  "

# Row 13
Set-TextValue $ws.Cells.Item(13, 4) "0.07860"
Set-TextValue $ws.Cells.Item(13, 5) "  -3.70%  "

# Row 14
Set-TextValue $ws.Cells.Item(14, 4) "1.819.39"
Set-TextValue $ws.Cells.Item(14, 5) "  -2.78%  "

# Row 15
Set-TextValue $ws.Cells.Item(15, 4) "5.075"
Set-TextValue $ws.Cells.Item(15, 5) "  -4.39%  "

# Row 16
Set-TextValue $ws.Cells.Item(16, 4) "87.68"
Set-TextValue $ws.Cells.Item(16, 5) "  -5.61%  "

# Row 17
Set-TextValue $ws.Cells.Item(17, 5) "  +0.04%  "

# Row 18
Set-TextValue $ws.Cells.Item(18, 4) "14.08"
Set-TextValue $ws.Cells.Item(18, 5) "  -5.16%  "

# Row 19
Set-TextValue $ws.Cells.Item(19, 2) "Dai"
Set-TextValue $ws.Cells.Item(19, 3) "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Cells.Item(19, 4) "1.001"
Set-TextValue $ws.Cells.Item(19, 5) "  +0.02%  "

# Row 20
Set-TextValue $ws.Cells.Item(20, 2) "ShibaInu"
Set-TextValue $ws.Cells.Item(20, 3) "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue $ws.Cells.Item(20, 4) "0.000008042"
Set-TextValue $ws.Cells.Item(20, 5) "  -5.45%  "

# Row 21
Set-TextValue $ws.Cells.Item(21, 4) "25.774.54"
Set-TextValue $ws.Cells.Item(21, 5) "  -3.95%  "

# Row 22
Set-TextValue $ws.Cells.Item(22, 4) "4.748"
Set-TextValue $ws.Cells.Item(22, 5) "  -4.84%  "

# Row 23
Set-TextValue $ws.Cells.Item(23, 5) "  -5.85%  "

# Row 24
Set-TextValue $ws.Cells.Item(24, 4) "6.109"
Set-TextValue $ws.Cells.Item(24, 5) "  -4.32%  "

# Row 25
Set-TextValue $ws.Cells.Item(25, 4) "2.239"
Set-TextValue $ws.Cells.Item(25, 5) "  -3.20%  "

# Row 26
Set-TextValue $ws.Cells.Item(26, 4) "142.38"
Set-TextValue $ws.Cells.Item(26, 5) "  -2.50%  "

# Row 27
Set-TextValue $ws.Cells.Item(27, 5) "  -3.92%  "

# Row 28
Set-TextValue $ws.Cells.Item(28, 4) "17.14"
Set-TextValue $ws.Cells.Item(28, 5) "  -4.96%  "

# Row 29
Set-TextValue $ws.Cells.Item(29, 4) "109.16"
Set-TextValue $ws.Cells.Item(29, 5) "  -4.27%  "

# Row 30
Set-TextValue $ws.Cells.Item(30, 4) "4.329"
Set-TextValue $ws.Cells.Item(30, 5) "  -8.32%  "

# Row 31
Set-TextValue $ws.Cells.Item(31, 4) "4.226"
Set-TextValue $ws.Cells.Item(31, 5) "  -8.75%  "

# Row 32
Set-TextValue $ws.Cells.Item(32, 4) "0.08752"
Set-TextValue $ws.Cells.Item(32, 5) "  -4.38%  "

# Row 33
Set-TextValue $ws.Cells.Item(33, 4) "0.04885"
Set-TextValue $ws.Cells.Item(33, 5) "  -2.88%  "

# Row 34
Set-TextValue $ws.Cells.Item(34, 4) "0.7283"
Set-TextValue $ws.Cells.Item(34, 5) "  -9.45%  "

# Row 35
Set-TextValue $ws.Cells.Item(35, 4) "1.135"
Set-TextValue $ws.Cells.Item(35, 5) "  -3.30%  "

# Row 36
Set-TextValue $ws.Cells.Item(36, 4) "2.884"
Set-TextValue $ws.Cells.Item(36, 5) "  -2.09%  "

# Row 37
Set-TextValue $ws.Cells.Item(37, 4) "3.146"
Set-TextValue $ws.Cells.Item(37, 5) "  -1.55%  "

# Row 38
Set-TextValue $ws.Cells.Item(38, 4) "2.357"
Set-TextValue $ws.Cells.Item(38, 5) "  -12.81%  "

# Row 39
Set-TextValue $ws.Cells.Item(39, 4) "0.01854"
Set-TextValue $ws.Cells.Item(39, 5) "  -5.06%  "

# Row 40
Set-TextValue $ws.Cells.Item(40, 5) "  -15.99%  "

# Row 41
Set-TextValue $ws.Cells.Item(41, 4) "0.9706"
Set-TextValue $ws.Cells.Item(41, 5) "  -8.93%  "

# Row 42
Set-TextValue $ws.Cells.Item(42, 4) "114.15"
Set-TextValue $ws.Cells.Item(42, 5) "  -0.58%  "

# Row 43
Set-TextValue $ws.Cells.Item(43, 4) "6.234"
Set-TextValue $ws.Cells.Item(43, 5) "  -4.51%  "

# Row 44
Set-TextValue $ws.Cells.Item(44, 4) "8.002"
Set-TextValue $ws.Cells.Item(44, 5) "  -8.99%  "

# Row 45
Set-TextValue $ws.Cells.Item(45, 5) "  +0.00%  "

# Row 46
Set-TextValue $ws.Cells.Item(46, 4) "0.4530"
Set-TextValue $ws.Cells.Item(46, 5) "  -13.53%  "

# Row 47
Set-TextValue $ws.Cells.Item(47, 4) "0.1367"
Set-TextValue $ws.Cells.Item(47, 5) "  -8.39%  "

# Row 48
Set-TextValue $ws.Cells.Item(48, 4) "36.40"
Set-TextValue $ws.Cells.Item(48, 5) "  -3.26%  "

# Row 49
Set-TextValue $ws.Cells.Item(49, 4) "9.201"
Set-TextValue $ws.Cells.Item(49, 5) "  -7.45%  "

# Row 50
Set-TextValue $ws.Cells.Item(50, 5) "  -9.12%  "

# Row 51
Set-TextValue $ws.Cells.Item(51, 4) "0.05832"
Set-TextValue $ws.Cells.Item(51, 5) "  -3.70%  "

